# Apply odds updates to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 3.1
$ws.Range("I2").Value = 2.5
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("AC2").Value = 7.5
$ws.Range("AE2").Value = 12
$ws.Range("AF2").Value = 34
$ws.Range("AI2").Value = 6.5
$ws.Range("AP2").Value = 10
$ws.Range("AQ2").Value = 23

# Row 3
$ws.Range("G3").Value = 1.6
$ws.Range("H3").Value = 3.9
$ws.Range("K3").Value = 2.25
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.42
$ws.Range("R3").Value = 2.9
$ws.Range("S3").Value = 1.85
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 2.48
$ws.Range("V3").Value = 1.54
$ws.Range("W3").Value = 3
$ws.Range("X3").Value = 1.36
$ws.Range("Y3").Value = 1.36
$ws.Range("Z3").Value = 3
$ws.Range("AA3").Value = 1.91
$ws.Range("AB3").Value = 1.91
$ws.Range("AC3").Value = 7
$ws.Range("AI3").Value = 11
$ws.Range("AJ3").Value = 7.5

# Row 4
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 3.7
$ws.Range("U4").Value = 4.4
$ws.Range("V4").Value = 1.21
$ws.Range("Y4").Value = 1.62
$ws.Range("Z4").Value = 2.2
$ws.Range("AA4").Value = 2.25
$ws.Range("AB4").Value = 1.57

# Row 11
$ws.Range("H11").Value = 2.88
